$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 172765
$ws.Range("E2").Value = 8292
$ws.Range("F2").Value = 8292
$ws.Range("G2").Value = 6730
$ws.Range("H2").Value = 4793
$ws.Range("I2").Value = 3313
$ws.Range("J2").Value = 1480
$ws.Range("K2").Value = 182446
$ws.Range("L2").Value = 114573
$ws.Range("M2").Value = 67873
$ws.Range("N2").Value = 52915
$ws.Range("O2").Value = 14958
$ws.Range("P2").Value = 5573
$ws.Range("Q2").Value = 4260
$ws.Range("R2").Value = 3246
$ws.Range("S2").Value = -802
$ws.Range("T2").Value = 1881
$ws.Range("U2").Value = 2380
$ws.Range("V2").Value = 25974
$ws.Range("W2").Value = 4.8
$ws.Range("X2").Value = 2.77
$ws.Range("Y2").Value = 6.55
$ws.Range("Z2").Value = 2.91
$ws.Range("AA2").Value = 168.81
$ws.Range("AB2").Value = 848.26
$ws.Range("AC2").Value = 2972
$ws.Range("AD2").Value = 14.16
$ws.Range("AE2").Value = 47477
$ws.Range("AF2").Value = 0.89
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.19
$ws.Range("AI2").Value = 16.82
$ws.Range("AJ2").Value = 111355765
$ws.Range("D3").Value = 192332
$ws.Range("E3").Value = 10893
$ws.Range("F3").Value = 10893
$ws.Range("G3").Value = 8981
$ws.Range("H3").Value = 6638
$ws.Range("I3").Value = 4210
$ws.Range("J3").Value = 2428
$ws.Range("K3").Value = 193501
$ws.Range("L3").Value = 119595
$ws.Range("M3").Value = 73906
$ws.Range("N3").Value = 57203
$ws.Range("O3").Value = 16703
$ws.Range("P3").Value = 5573
$ws.Range("Q3").Value = 5960
$ws.Range("R3").Value = -11738
$ws.Range("S3").Value = 307
$ws.Range("T3").Value = 1865
$ws.Range("U3").Value = 4095
$ws.Range("V3").Value = 26327
$ws.Range("W3").Value = 5.66
$ws.Range("X3").Value = 3.45
$ws.Range("Y3").Value = 7.65
$ws.Range("Z3").Value = 3.53
$ws.Range("AA3").Value = 161.82
$ws.Range("AB3").Value = 919.99
$ws.Range("AC3").Value = 3777
$ws.Range("AD3").Value = 7.56
$ws.Range("AE3").Value = 51324
$ws.Range("AF3").Value = 0.5600000000000001
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.75
$ws.Range("AI3").Value = 13.24
$ws.Range("AJ3").Value = 111355765
$ws.Range("D4").Value = 188250
$ws.Range("E4").Value = 11590
$ws.Range("F4").Value = 11590
$ws.Range("G4").Value = 9815
$ws.Range("H4").Value = 7317
$ws.Range("I4").Value = 5721
$ws.Range("J4").Value = 1596
$ws.Range("K4").Value = 198734
$ws.Range("L4").Value = 117488
$ws.Range("M4").Value = 81246
$ws.Range("N4").Value = 63290
$ws.Range("O4").Value = 17956
$ws.Range("P4").Value = 5573
$ws.Range("Q4").Value = 10865
$ws.Range("R4").Value = -7752
$ws.Range("S4").Value = -1769
$ws.Range("T4").Value = 1787
$ws.Range("U4").Value = 9079
$ws.Range("V4").Value = 25599
$ws.Range("W4").Value = 6.16
$ws.Range("X4").Value = 3.89
$ws.Range("Y4").Value = 9.5
$ws.Range("Z4").Value = 3.73
$ws.Range("AA4").Value = 144.61
$ws.Range("AB4").Value = 1015.93
$ws.Range("AC4").Value = 5133
$ws.Range("AD4").Value = 8.34
$ws.Range("AE4").Value = 56786
$ws.Range("AF4").Value = 0.75
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 1.17
$ws.Range("AI4").Value = 9.74
$ws.Range("AJ4").Value = 111355765
$ws.Range("D5").Value = 168871
$ws.Range("E5").Value = 9861
$ws.Range("F5").Value = 9861
$ws.Range("G5").Value = 5548
$ws.Range("H5").Value = 3716
$ws.Range("I5").Value = 2017
$ws.Range("J5").Value = 1699
$ws.Range("K5").Value = 184319
$ws.Range("L5").Value = 99590
$ws.Range("M5").Value = 84729
$ws.Range("N5").Value = 64676
$ws.Range("O5").Value = 20053
$ws.Range("P5").Value = 5573
$ws.Range("Q5").Value = 5144
$ws.Range("R5").Value = -163
$ws.Range("S5").Value = -3809
$ws.Range("T5").Value = 498
$ws.Range("U5").Value = 4646
$ws.Range("V5").Value = 22818
$ws.Range("W5").Value = 5.84
$ws.Range("X5").Value = 2.2
$ws.Range("Y5").Value = 3.15
$ws.Range("Z5").Value = 1.94
$ws.Range("AA5").Value = 117.54
$ws.Range("AB5").Value = 1048.98
$ws.Range("AC5").Value = 1810
$ws.Range("AD5").Value = 20.06
$ws.Range("AE5").Value = 58029
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 1.38
$ws.Range("AI5").Value = 27.63
$ws.Range("AJ5").Value = 111355765
$ws.Range("D6").Value = 167309
$ws.Range("E6").Value = 8400
$ws.Range("F6").Value = 8400
$ws.Range("G6").Value = 8979
$ws.Range("H6").Value = 5353
$ws.Range("I6").Value = 3816
$ws.Range("K6").Value = 180546
$ws.Range("L6").Value = 97628
$ws.Range("M6").Value = 82919
$ws.Range("N6").Value = 62741
$ws.Range("P6").Value = 5573
$ws.Range("Q6").Value = 2495
$ws.Range("R6").Value = -920
$ws.Range("S6").Value = -373
$ws.Range("T6").Value = 1919
$ws.Range("U6").Value = 576
$ws.Range("V6").Value = 24067
$ws.Range("W6").Value = 5.02
$ws.Range("X6").Value = 3.2
$ws.Range("Y6").Value = 5.99
$ws.Range("Z6").Value = 2.93
$ws.Range("AA6").Value = 117.74
$ws.Range("AB6").Value = 1074.89
$ws.Range("AC6").Value = 3424
$ws.Range("AD6").Value = 15.95
$ws.Range("AE6").Value = 56293
$ws.Range("AF6").Value = 0.97
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 0.92
$ws.Range("AI6").Value = 14.6
$ws.Range("AJ6").Value = 111355765
$ws.Range("D7").Value = 170502
$ws.Range("E7").Value = 9080
$ws.Range("G7").Value = 9562
$ws.Range("H7").Value = 7117
$ws.Range("I7").Value = 5199
$ws.Range("K7").Value = 186301
$ws.Range("L7").Value = 96654
$ws.Range("M7").Value = 89647
$ws.Range("N7").Value = 68248
$ws.Range("P7").Value = 5571
$ws.Range("Q7").Value = 7090
$ws.Range("R7").Value = -569
$ws.Range("S7").Value = -1749
$ws.Range("T7").Value = 1265
$ws.Range("U7").Value = 5258
$ws.Range("W7").Value = 5.33
$ws.Range("X7").Value = 4.17
$ws.Range("Y7").Value = 7.94
$ws.Range("Z7").Value = 3.88
$ws.Range("AA7").Value = 107.82
$ws.Range("AC7").Value = 4665
$ws.Range("AD7").Value = 9.01
$ws.Range("AE7").Value = 61234
$ws.Range("AF7").Value = 0.6899999999999999
$ws.Range("AG7").Value = 530
$ws.Range("AH7").Value = 1.26
$ws.Range("AI7").Value = 11.34
$ws.Range("D8").Value = 176219
$ws.Range("E8").Value = 10064
$ws.Range("G8").Value = 9974
$ws.Range("H8").Value = 7311
$ws.Range("I8").Value = 5356
$ws.Range("K8").Value = 191950
$ws.Range("L8").Value = 97026
$ws.Range("M8").Value = 94925
$ws.Range("N8").Value = 72357
$ws.Range("P8").Value = 5571
$ws.Range("Q8").Value = 7207
$ws.Range("R8").Value = -2894
$ws.Range("S8").Value = -1835
$ws.Range("T8").Value = 1426
$ws.Range("U8").Value = 4899
$ws.Range("W8").Value = 5.71
$ws.Range("X8").Value = 4.15
$ws.Range("Y8").Value = 7.62
$ws.Range("Z8").Value = 3.86
$ws.Range("AA8").Value = 102.21
$ws.Range("AC8").Value = 4806
$ws.Range("AD8").Value = 7.95
$ws.Range("AE8").Value = 64921
$ws.Range("AF8").Value = 0.59
$ws.Range("AG8").Value = 574
$ws.Range("AH8").Value = 1.5
$ws.Range("AI8").Value = 11.93
$ws.Range("D9").Value = 185428
$ws.Range("E9").Value = 10822
$ws.Range("G9").Value = 10832
$ws.Range("H9").Value = 7953
$ws.Range("I9").Value = 5898
$ws.Range("K9").Value = 199679
$ws.Range("L9").Value = 97961
$ws.Range("M9").Value = 101634
$ws.Range("N9").Value = 77638
$ws.Range("P9").Value = 5571
$ws.Range("Q9").Value = 6962
$ws.Range("R9").Value = -3121
$ws.Range("S9").Value = -1588
$ws.Range("T9").Value = 1453
$ws.Range("U9").Value = 4850
$ws.Range("W9").Value = 5.84
$ws.Range("X9").Value = 4.29
$ws.Range("Y9").Value = 7.86
$ws.Range("Z9").Value = 4.06
$ws.Range("AA9").Value = 96.39
$ws.Range("AC9").Value = 5292
$ws.Range("AD9").Value = 7.22
$ws.Range("AE9").Value = 69659
$ws.Range("AF9").Value = 0.55
$ws.Range("AG9").Value = 610
$ws.Range("AH9").Value = 1.6
$ws.Range("AI9").Value = 11.51
